# Generate Report for Handback
# Updates the localization-status workbook after a handback event:
#  - Overview sheet: zh-cn / de-de status columns move from "Ready for
#    handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets: Status updated the same way, the
#    "Latest Handback DateTime" is refreshed to the new handback
#    timestamp, and the (now resolved) Error Detail is cleared.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-13 13:02:23"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-13 13:02:33"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- Overview sheet column widths (zh-cn / de-de status columns widen
#      to fit the longer "Handed back..." text) ----
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668
